$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151 (shifts existing rows 151.. down by one).
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new record.
$ws.Cells.Item(151, 1).Value = 4
$ws.Cells.Item(151, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(151, 3).Value = 'Los Lagos'
$ws.Cells.Item(151, 4).Value = 44522
$ws.Cells.Item(151, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(151, 5).Value = 10
$ws.Cells.Item(151, 6).Value = 100112003
$ws.Cells.Item(151, 7).Value = 'Ajo'
$ws.Cells.Item(151, 8).Value = 'Chino'
$ws.Cells.Item(151, 9).Value = 'Primera'
$ws.Cells.Item(151, 10).Value = 80
$ws.Cells.Item(151, 11).Value = 21000
$ws.Cells.Item(151, 12).Value = 22000
$ws.Cells.Item(151, 13).Value = 21500
$ws.Cells.Item(151, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(151, 15).Value = 'China'
$ws.Cells.Item(151, 16).Value = 2150
$ws.Cells.Item(151, 17).Value = 10
$ws.Cells.Item(151, 18).Value = 'Hortaliza'
